# Generate Report for Handback
#
# - Overview / zh-cn / de-de sheets: status text "Ready for handoff" ->
#   "Handed back: in sync with en-US"
# - zh-cn sheet: refresh the "Latest Handback DateTime" timestamps
# - de-de sheet: the handback for this round is now in sync with en-US, so
#   it also gets a "Latest Target File" / "Latest Handback File" pair of
#   hyperlinked entries (mirroring the zh-cn sheet's layout) plus a refreshed
#   "Latest Handback DateTime" timestamp.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---- Overview sheet ------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# ---- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = $statusHandedBack
$zhcn.Range("G2").Value = "2016-03-09 02:43:58"
$zhcn.Range("B3").Value = $statusHandedBack
$zhcn.Range("G3").Value = "2016-03-09 02:43:58"

# ---- de-de sheet -----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = $statusHandedBack
$dede.Range("B3").Value = $statusHandedBack

$targetFileName = "a.md"
$handbackFileName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$targetUrl2 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1f1d39b0c1ce2db6f3bb7bba3f0c0dcbb2a8e9b1/e2e/a.md"
$handbackUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3e7b1f0a4c9d2e8f5a6b1c0d9e8f7a6b5c4d3e2f/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$targetUrl3 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1f1d39b0c1ce2db6f3bb7bba3f0c0dcbb2a8e9b1/e2e/a.md"
$handbackUrl3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3e7b1f0a4c9d2e8f5a6b1c0d9e8f7a6b5c4d3e2f/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("E2"), $targetUrl2, "", "", $targetFileName)
$dede.Hyperlinks.Add($dede.Range("F2"), $handbackUrl2, "", "", $handbackFileName)
$dede.Range("G2").Value = "2016-03-09 02:44:33"

$dede.Hyperlinks.Add($dede.Range("E3"), $targetUrl3, "", "", $targetFileName)
$dede.Hyperlinks.Add($dede.Range("F3"), $handbackUrl3, "", "", $handbackFileName)
$dede.Range("G3").Value = "2016-03-09 02:44:33"
